$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric-looking "Price" values in column D (force text to preserve exact formatting) ---
$priceUpdates = @{
    'D2' = '243.54'
    'D3' = '23.77'
    'D4' = '5.315'
    'D5' = '0.05771'
    'D6' = '6.479'
    'D7' = '3.347'
    'D8' = '0.8113'
    'D9' = '0.8876'
    'D11' = '0.07358'
    'D12' = '0.03087'
    'D13' = '0.03067'
    'D14' = '0.09333'
    'D15' = '3.877'
    'D16' = '0.001542'
    'D17' = '0.04697'
    'D18' = '0.0006020'
    'D19' = '0.006182'
    'D20' = '0.001296'
    'D21' = '0.00008818'
    'D22' = '3.582'
    'D28' = '0.0002354'
    'D40' = '0.03765'
    'D41' = '0.006383'
    'D42' = '0.1051'
    'D43' = '0.003289'
    'D44' = '0.007607'
    'D45' = '0.00005481'
    'D46' = '0.00000000752'
    'D47' = '0.5509'
    'D48' = '0.001850'
    'D49' = '0.00002104'
    'D50' = '0.0002004'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# --- Update plain text cells (Coin name, Link, Volume label) ---
$textUpdates = @{
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'E41' = '40KickTokenKICK'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'E42' = '41BKEXTokenBKK'
    'B43' = 'CEJI'
    'C43' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'E43' = '42CEJICEJI'
    'E48' = '47BOLOBOLOWorstin24h'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}
